# Weekly update for Fruta / Piña - Vega Central Mapocho de Santiago.
# A new week of data (2021-11-04, Caramelo, Especial/Primera/Segunda/Tercera)
# is inserted right before the existing 2021-07-07 block (row 664), pushing
# every following row down by 4. The four rows that fall off the bottom of
# the old range keep their original values at the new bottom rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows at 664..667; everything below shifts down by 4.
$ws.Rows("664:667").Insert()

# Row 664 - Especial
$ws.Range("A664").Value = 9
$ws.Range("B664").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C664").Value = "Metropolitana"
$ws.Range("D664").Value = 44504
$ws.Range("E664").Value = 13
$ws.Range("F664").Value = "Fruta"
$ws.Range("G664").Value = 100108
$ws.Range("H664").Value = "Tropicales y subtropicales"
$ws.Range("I664").Value = 100108005
$ws.Range("J664").Value = "Piña"
$ws.Range("K664").Value = "Caramelo"
$ws.Range("L664").Value = "Especial"
$ws.Range("M664").Value = 25
$ws.Range("N664").Value = 18000
$ws.Range("O664").Value = 18000
$ws.Range("P664").Value = 18000
$ws.Range("Q664").Value = "$/caja 10 unidades"
$ws.Range("R664").Value = "Ecuador"
$ws.Range("S664").Value = 1800
$ws.Range("T664").Value = 10

# Row 665 - Primera
$ws.Range("A665").Value = 9
$ws.Range("B665").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C665").Value = "Metropolitana"
$ws.Range("D665").Value = 44504
$ws.Range("E665").Value = 13
$ws.Range("F665").Value = "Fruta"
$ws.Range("G665").Value = 100108
$ws.Range("H665").Value = "Tropicales y subtropicales"
$ws.Range("I665").Value = 100108005
$ws.Range("J665").Value = "Piña"
$ws.Range("K665").Value = "Caramelo"
$ws.Range("L665").Value = "Primera"
$ws.Range("M665").Value = 20
$ws.Range("N665").Value = 18000
$ws.Range("O665").Value = 18000
$ws.Range("P665").Value = 18000
$ws.Range("Q665").Value = "$/caja 12 unidades"
$ws.Range("R665").Value = "Ecuador"
$ws.Range("S665").Value = 1500
$ws.Range("T665").Value = 12

# Row 666 - Segunda
$ws.Range("A666").Value = 9
$ws.Range("B666").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C666").Value = "Metropolitana"
$ws.Range("D666").Value = 44504
$ws.Range("E666").Value = 13
$ws.Range("F666").Value = "Fruta"
$ws.Range("G666").Value = 100108
$ws.Range("H666").Value = "Tropicales y subtropicales"
$ws.Range("I666").Value = 100108005
$ws.Range("J666").Value = "Piña"
$ws.Range("K666").Value = "Caramelo"
$ws.Range("L666").Value = "Segunda"
$ws.Range("M666").Value = 35
$ws.Range("N666").Value = 18000
$ws.Range("O666").Value = 18000
$ws.Range("P666").Value = 18000
$ws.Range("Q666").Value = "$/caja 14 unidades"
$ws.Range("R666").Value = "Ecuador"
$ws.Range("S666").Value = 1286
$ws.Range("T666").Value = 14

# Row 667 - Tercera
$ws.Range("A667").Value = 9
$ws.Range("B667").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C667").Value = "Metropolitana"
$ws.Range("D667").Value = 44504
$ws.Range("E667").Value = 13
$ws.Range("F667").Value = "Fruta"
$ws.Range("G667").Value = 100108
$ws.Range("H667").Value = "Tropicales y subtropicales"
$ws.Range("I667").Value = 100108005
$ws.Range("J667").Value = "Piña"
$ws.Range("K667").Value = "Caramelo"
$ws.Range("L667").Value = "Tercera"
$ws.Range("M667").Value = 25
$ws.Range("N667").Value = 18000
$ws.Range("O667").Value = 18000
$ws.Range("P667").Value = 18000
$ws.Range("Q667").Value = "$/caja 16 unidades"
$ws.Range("R667").Value = "Ecuador"
$ws.Range("S667").Value = 1125
$ws.Range("T667").Value = 16
